# New crime data collected — update the weekly CompStat figures for the
# 123rd Precinct report: header volume/date strings, the crime-stat table
# (rows 16-27) and the bestFit width of column E that shifts as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header strings: "Volume 30 Number 5" -> "... Number 6" and the covered
# week dates move forward one week.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 30   Number  6"
$ws.Range("C9").Value  = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# ---------------------------------------------------------------------------
# Crime-complaint grid (rows 16-27). Most cells simply get new numbers;
# a few cells flip between a numeric value and the literal placeholder text
# ("0" / "***.*") used elsewhere in the sheet for "no data" rows, so those
# need the cell pre-formatted as Text before the value is written (otherwise
# Excel auto-coerces "0" back into a number).
# ---------------------------------------------------------------------------

# Row 16
$ws.Range("N16").Value = -80

# Row 17
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 3
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -22.222222222222
$ws.Range("L17").Value = 133.333333333333
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 40

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 350
$ws.Range("L18").Value = 28.571428571428
$ws.Range("M18").Value = -18.181818181818
$ws.Range("N18").Value = -76.315789473684

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 5.555555555555
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = 31
$ws.Range("K19").Value = -6.451612903225
$ws.Range("L19").Value = 61.111111111111
$ws.Range("M19").Value = 38.095238095238
$ws.Range("N19").Value = 107.142857142857

# Row 20 — C20 goes from a number to the text placeholder "0"; D20/E20 go
# the other way, from text placeholders to real numbers.
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 2
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -50
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = -9.090909090909
$ws.Range("N20").Value = -88.372093023255

# Row 21
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -41.666666666666
$ws.Range("F21").Value = 29
$ws.Range("H21").Value = -9.375
$ws.Range("I21").Value = 56
$ws.Range("J21").Value = 53
$ws.Range("K21").Value = 5.660377358490
$ws.Range("M21").Value = 21.739130434782
$ws.Range("N21").Value = -62.162162162162

# Row 24
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -43.75
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = 11.111111111111
$ws.Range("I24").Value = 66
$ws.Range("J24").Value = 42
$ws.Range("K24").Value = 57.142857142857
$ws.Range("L24").Value = 175
$ws.Range("M24").Value = -2.941176470588

# Row 25
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 87.5
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 19
$ws.Range("K25").Value = 5.263157894736
$ws.Range("L25").Value = 566.666666666667
$ws.Range("M25").Value = -16.666666666666

# Row 27 — inverse of row 20: C27 goes from text placeholder "0" to a real
# number; D27/E27 go from real numbers to the text placeholders.
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = 100

# ---------------------------------------------------------------------------
# Column E narrows slightly (bestFit) now that its numbers are shorter.
# ColumnWidth is quantized to whole pixels by the host, so this is the
# closest achievable value to the authored 7.433768 "width" units.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.71
